# Auto-generated Excel COM-interop script
# Applies cell-value updates to the Lamia_Profits profit-tracking workbook
# (computed market-price / profit columns H:N refreshed by the scheduled pricing runner)

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 9
$ws.Range("H9").Value = 185.875
$ws.Range("J9").Value = 144.75
$ws.Range("L9").Value = 144.75
$ws.Range("N9").Value = -482.75
# Row 21
$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").ClearContents()
$ws.Range("N21").Value = 0
# Row 23
$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("L23").ClearContents()
$ws.Range("N23").Value = 0
# Row 29
$ws.Range("H29").Value = 5463.5
$ws.Range("J29").Value = 6429.143
$ws.Range("L29").Value = 19287.429
$ws.Range("N29").Value = -19849.429
# Row 38
$ws.Range("H38").Value = 948.6667
$ws.Range("I38").Value = 38.4
$ws.Range("K38").Value = 115.2
$ws.Range("M38").Value = 256.8
# Row 43
$ws.Range("H43").Value = 5632.1875
$ws.Range("I43").Value = 3096.25
$ws.Range("K43").Value = 3096.25
$ws.Range("M43").Value = -3027.25
# Row 112
$ws.Range("H112").Value = 1400.8948
$ws.Range("I112").Value = 935.75
$ws.Range("J112").Value = 1524.9333
$ws.Range("K112").Value = 2807.25
$ws.Range("L112").Value = 4574.7999
$ws.Range("M112").Value = -1699.25
$ws.Range("N112").Value = -6790.7999
# Row 137
$ws.Range("H137").Value = 3090.9385
$ws.Range("I137").Value = 2437.303
$ws.Range("J137").Value = 3765
$ws.Range("K137").Value = 7311.909
$ws.Range("L137").Value = 11295
$ws.Range("M137").Value = -4761.909
$ws.Range("N137").Value = -16395

$ws = $wb.Worksheets.Item("ARM")
# Row 38
$ws.Range("H38").Value = 8007.6
$ws.Range("I38").Value = 8007.6
$ws.Range("K38").Value = 8007.6
$ws.Range("M38").Value = -7540.6
# Row 45
$ws.Range("H45").Value = 90914140
$ws.Range("I45").Value = 125003080
$ws.Range("J45").Value = 10309.333
$ws.Range("K45").Value = 125003080
$ws.Range("L45").Value = 10309.333
$ws.Range("M45").Value = -125002703
$ws.Range("N45").Value = -11063.333
# Row 74
$ws.Range("H74").Value = 8550337
$ws.Range("I74").Value = 10755359
$ws.Range("J74").Value = 5876
$ws.Range("K74").Value = 10755359
$ws.Range("L74").Value = 5876
$ws.Range("M74").Value = -10754485
$ws.Range("N74").Value = -7624
# Row 77
$ws.Range("H77").Value = 8550337
$ws.Range("I77").Value = 10755359
$ws.Range("J77").Value = 5876
$ws.Range("K77").Value = 53776795
$ws.Range("L77").Value = 29380
$ws.Range("M77").Value = -53772427
$ws.Range("N77").Value = -38116
# Row 97
$ws.Range("H97").Value = 3976.72
$ws.Range("I97").Value = 3117
$ws.Range("K97").Value = 3117
$ws.Range("M97").Value = -2621
# Row 122
$ws.Range("H122").Value = 4098.826
$ws.Range("I122").Value = 3721.5334
$ws.Range("J122").Value = 4806.25
$ws.Range("K122").Value = 11164.6002
$ws.Range("L122").Value = 14418.75
$ws.Range("M122").Value = -8714.600199999999
$ws.Range("N122").Value = -19318.75
# Row 132
$ws.Range("H132").Value = 4044.9822
$ws.Range("I132").Value = 2913.1345
$ws.Range("J132").Value = 18759
$ws.Range("K132").Value = 8739.4035
$ws.Range("L132").Value = 56277
$ws.Range("M132").Value = -6209.4035
$ws.Range("N132").Value = -61337

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 2746
$ws.Range("I20").Value = 2392.111
$ws.Range("J20").Value = 4338.5
$ws.Range("K20").Value = 2392.111
$ws.Range("L20").Value = 4338.5
$ws.Range("M20").Value = -2145.111
$ws.Range("N20").Value = -4832.5
# Row 22
$ws.Range("H22").Value = 1556.5
$ws.Range("I22").Value = 1008.375
$ws.Range("J22").Value = 3749
$ws.Range("K22").Value = 1008.375
$ws.Range("L22").Value = 3749
$ws.Range("M22").Value = -835.375
$ws.Range("N22").Value = -4095
# Row 29
$ws.Range("H29").Value = 4300.3335
$ws.Range("I29").Value = 4001
$ws.Range("K29").Value = 4001
$ws.Range("M29").Value = -3712
# Row 94
$ws.Range("H94").Value = 657.8889
$ws.Range("I94").Value = 657.8889
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 657.8889
$ws.Range("L94").Value = 0
$ws.Range("M94").ClearContents()
$ws.Range("N94").Value = -206.8889
# Row 126
$ws.Range("H126").Value = 77800
$ws.Range("J126").Value = 77800
$ws.Range("L126").Value = 77800
$ws.Range("N126").Value = -87680

$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 224.875
$ws.Range("I7").Value = 212
$ws.Range("J7").Value = 237.75
$ws.Range("K7").Value = 212
$ws.Range("L7").Value = 237.75
$ws.Range("M7").Value = -99
$ws.Range("N7").Value = -463.75
# Row 22
$ws.Range("H22").Value = 1514.9166
$ws.Range("I22").Value = 656
$ws.Range("J22").Value = 4091.6667
$ws.Range("K22").Value = 656
$ws.Range("L22").Value = 4091.6667
$ws.Range("M22").Value = -306
$ws.Range("N22").Value = -4791.6667
# Row 58
$ws.Range("H58").Value = 4141.8423
$ws.Range("I58").Value = 1160.0769
$ws.Range("J58").Value = 10602.333
$ws.Range("K58").Value = 1160.0769
$ws.Range("L58").Value = 10602.333
$ws.Range("M58").Value = -957.0769
$ws.Range("N58").Value = -11008.333
# Row 86
$ws.Range("H86").Value = 7145.2
$ws.Range("I86").Value = 3750
$ws.Range("J86").Value = 12238
$ws.Range("K86").Value = 3750
$ws.Range("L86").Value = 12238
$ws.Range("M86").Value = -2627
$ws.Range("N86").Value = -14484
# Row 89
$ws.Range("H89").Value = 7145.2
$ws.Range("I89").Value = 3750
$ws.Range("J89").Value = 12238
$ws.Range("K89").Value = 18750
$ws.Range("L89").Value = 61190
$ws.Range("M89").Value = -13134
$ws.Range("N89").Value = -72422
# Row 132
$ws.Range("H132").Value = 2656.9
$ws.Range("I132").Value = 2109.7837
$ws.Range("J132").Value = 9404.666999999999
$ws.Range("K132").Value = 6329.3511
$ws.Range("L132").Value = 28214.001
$ws.Range("M132").Value = -3799.3511
$ws.Range("N132").Value = -33274.001
# Row 136
$ws.Range("H136").Value = 4141.8423
$ws.Range("I136").Value = 1160.0769
$ws.Range("J136").Value = 10602.333
$ws.Range("K136").Value = 3480.2307
$ws.Range("L136").Value = 31806.999
$ws.Range("M136").Value = -930.2307000000001
$ws.Range("N136").Value = -36906.999

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 6574.522
$ws.Range("I5").Value = 614.6
$ws.Range("J5").Value = 17749.375
$ws.Range("K5").Value = 1843.8
$ws.Range("L5").Value = 53248.125
$ws.Range("M5").Value = -1731.8
$ws.Range("N5").Value = -53472.125
# Row 113
$ws.Range("H113").Value = 1341.2069
$ws.Range("I113").Value = 934.9
$ws.Range("J113").Value = 1555.0526
$ws.Range("K113").Value = 2804.7
$ws.Range("L113").Value = 4665.1578
$ws.Range("M113").Value = -634.6999999999998
$ws.Range("N113").Value = -9005.157800000001
# Row 132
$ws.Range("H132").Value = 6719.25
$ws.Range("I132").Value = 5500
$ws.Range("K132").Value = 49500
$ws.Range("M132").Value = -46970
# Row 133
$ws.Range("H133").Value = 23000
$ws.Range("J133").Value = 25000
$ws.Range("L133").Value = 75000
$ws.Range("N133").Value = -85120
# Row 135
$ws.Range("H135").Value = 6574.522
$ws.Range("I135").Value = 614.6
$ws.Range("J135").Value = 17749.375
$ws.Range("K135").Value = 5531.400000000001
$ws.Range("L135").Value = 159744.375
$ws.Range("M135").Value = -2996.400000000001
$ws.Range("N135").Value = -164814.375

$ws = $wb.Worksheets.Item("GSM")
# Row 113
$ws.Range("H113").Value = 5486.4116
$ws.Range("I113").Value = 2326.9
$ws.Range("J113").Value = 10000
$ws.Range("K113").Value = 2326.9
$ws.Range("L113").Value = 10000
$ws.Range("M113").Value = -156.9000000000001
$ws.Range("N113").Value = -14340
# Row 126
$ws.Range("H126").Value = 6513.25
$ws.Range("I126").Value = 1821.2858
$ws.Range("J126").Value = 10162.556
$ws.Range("K126").Value = 5463.857400000001
$ws.Range("L126").Value = 30487.668
$ws.Range("M126").Value = -2993.857400000001
$ws.Range("N126").Value = -35427.66800000001

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 1293.1111
$ws.Range("I16").Value = 1186
$ws.Range("J16").Value = 2150
$ws.Range("K16").Value = 1186
$ws.Range("L16").Value = 2150
$ws.Range("M16").Value = -1016
$ws.Range("N16").Value = -2490
# Row 82
$ws.Range("H82").Value = 5097.3
$ws.Range("I82").Value = 4204.7
$ws.Range("J82").Value = 5989.9
$ws.Range("K82").Value = 4204.7
$ws.Range("L82").Value = 5989.9
$ws.Range("M82").Value = -3843.7
$ws.Range("N82").Value = -6711.9
# Row 85
$ws.Range("H85").Value = 5097.3
$ws.Range("I85").Value = 4204.7
$ws.Range("J85").Value = 5989.9
$ws.Range("K85").Value = 4204.7
$ws.Range("L85").Value = 5989.9
$ws.Range("M85").Value = -2956.7
$ws.Range("N85").Value = -8485.9
# Row 93
$ws.Range("H93").Value = 12036.763
$ws.Range("I93").Value = 10445.818
$ws.Range("J93").Value = 14224.3125
$ws.Range("K93").Value = 10445.818
$ws.Range("L93").Value = 14224.3125
$ws.Range("M93").Value = -9197.817999999999
$ws.Range("N93").Value = -16720.3125

$ws = $wb.Worksheets.Item("WVR")
# Row 113
$ws.Range("H113").Value = 1089.075
$ws.Range("I113").Value = 983.52
$ws.Range("K113").Value = 2950.56
$ws.Range("M113").Value = -780.5599999999999
# Row 136
$ws.Range("H136").Value = 4818.5586
$ws.Range("I136").Value = 4128
$ws.Range("J136").Value = 7062.875
$ws.Range("K136").Value = 12384
$ws.Range("L136").Value = 21188.625
$ws.Range("M136").Value = -9834
$ws.Range("N136").Value = -26288.625
